# Add a new slide "Waarom een bepaalde tool gebruiken?" as the 4th slide,
# using the existing "Titel en object" (Title and Content) layout — the
# same layout used by the other slides in this deck.

$p = $ppt.ActivePresentation

# Layout index 2 == ppt/slideLayouts/slideLayout2.xml ("Titel en object"),
# the same layout already used by slides 1-3.
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title placeholder ---------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Name = "Titel 1"
$title.TextFrame.TextRange.Text = "Waarom een bepaalde tool gebruiken?"

# --- Content placeholder --------------------------------------------------
$content = $s.Shapes.Item(2)
$content.Name = "Tijdelijke aanduiding voor inhoud 2"

$bullets = @(
    "Open-source",
    "Gebruiksvriendelijke User Interface",
    "Kostenefficiënt",
    "Support",
    "Actieve update cycle",
    " "
)
$content.TextFrame.TextRange.Text = [string]::Join([char]13, $bullets)
